$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

$ws.Cells.Item($row, 1).Value2 = 11
$ws.Cells.Item($row, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value2 = "Bíobío"
$ws.Cells.Item($row, 4).Value2 = 45191
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value2 = 8
$ws.Cells.Item($row, 6).Value2 = 100112031
$ws.Cells.Item($row, 7).Value2 = "Poroto verde"
$ws.Cells.Item($row, 8).Value2 = "Magnum"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 50
$ws.Cells.Item($row, 11).Value2 = 16000
$ws.Cells.Item($row, 12).Value2 = 16000
$ws.Cells.Item($row, 13).Value2 = 16000
$ws.Cells.Item($row, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item($row, 15).Value2 = "Perú"
$ws.Cells.Item($row, 16).Value2 = 640
$ws.Cells.Item($row, 17).Value2 = 25
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
